# Update to README with pointers to Github project page
#
# The closing paragraph ("You can find everything I showed here at: ...")
# was left as a plain, un-styled paragraph. Bring it in line with the
# other entries in the "Some additional pointers" bulleted list right
# above it: give it the "Compact" paragraph style and continue the same
# bulleted list (numId 7) those sibling paragraphs use.

$d = $word.ActiveDocument

$needle = "You can find everything I showed here at:"

# Locate the target paragraph and the list item immediately preceding it
# (whose numbering we want to continue) by walking Paragraphs - this is
# more reliable in this host than slicing a Range returned by Find.
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($paraText -eq $needle) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target paragraph '$needle'"
}

$target = $d.Paragraphs.Item($targetIndex)
$precedingListItem = $d.Paragraphs.Item($targetIndex - 1)

# Match the "Compact" style used by the rest of the pointer list.
$target.Style = "Compact"

# Continue the existing bulleted list (same numId) rather than starting
# a brand new one.
$target.Range.ListFormat.ApplyListTemplate($precedingListItem.Range.ListFormat.ListTemplate, $true)

Write-Output "Styled paragraph $targetIndex as Compact list item"
